# Update "想去人数" (want-to-go count) values in column F
# for the "展览" and "全部类型" worksheets, reflecting the
# refreshed data pulled from bilibili at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Row (as it appears on the "展览" sheet) -> new value for column F
$updatesZhanLan = @{
    3  = 558
    4  = 1123
    6  = 64
    8  = 51
    9  = 1155
    10 = 16185
    12 = 198
    14 = 6324
    15 = 637
    16 = 124
    21 = 21
    29 = 46
    30 = 5038
    32 = 11275
    37 = 3831
}

# "全部类型" contains one extra row (an additional event) inserted
# before row 32 of "展览", so rows from that point on are shifted by +1.
$updatesQuanBu = @{
    3  = 558
    4  = 1123
    6  = 64
    8  = 51
    9  = 1155
    10 = 16185
    12 = 198
    14 = 6324
    15 = 637
    16 = 124
    21 = 21
    29 = 46
    30 = 5038
    33 = 11275
    38 = 3831
}

$wsZhanLan = $wb.Worksheets.Item("展览")
foreach ($row in $updatesZhanLan.Keys) {
    $wsZhanLan.Cells.Item($row, 6).Value = $updatesZhanLan[$row]
}

$wsQuanBu = $wb.Worksheets.Item("全部类型")
foreach ($row in $updatesQuanBu.Keys) {
    $wsQuanBu.Cells.Item($row, 6).Value = $updatesQuanBu[$row]
}
